$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DATA_USER_FORM"

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "NAME"
$ws.Cells.Item(1,2).Value = "SURNAME"
$ws.Cells.Item(1,3).Value = "COMPANY"
$ws.Cells.Item(1,4).Value = "EMAIL"

# ---- Data rows (name, surname, company, email) ----
$rows = @(
    @("Antonio", "Silva", "Ultranet DataCenters Inc.", "tonin@coldmail.froz.en"),
    @("Claudiene", "Santos", "FreeSpeech multimedia", "leitte_anonima@lettersweb.com"),
    @("Zywywz", "Correia de Campos", "Z y Z Contabilidade S/A", "zivio.do.agape@jinnmail.comma"),
    @("Maria Rosângela", "Américo", "Mari Rosa Engenharia Civil LTDA", "roseamerico@jinnmail.com"),
    @("Wellington", "Pereira", "FunkAgenda Produções", "wbeatzz@coldmail.froz.en"),
    @("Lucas", "Motta", "Galeria22 Exposições e Curadoria LTDA", "lucas_motta_21@ibapo.ru"),
    @("Karim", "Gharib", "Restaurante Al-Bait Maghribiy LTDA", "krmghrib_morroco@jinnmail.comma")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    $ws.Cells.Item($r,1).Value = $rowVals[0]
    $ws.Cells.Item($r,2).Value = $rowVals[1]
    $ws.Cells.Item($r,3).Value = $rowVals[2]
    $ws.Cells.Item($r,4).Value = $rowVals[3]
}

# Row 9 - Bruno Kunrath
$ws.Cells.Item(9,1).Value = "Bruno"
$ws.Cells.Item(9,2).Value = "Kunrath"
$ws.Cells.Item(9,3).Value = "Ultranet DataCenters Inc."
$ws.Cells.Item(9,4).Value = "bruko_ultra@coldmail.froz.en"

# Row 10 - Claudiene (no surname)
$ws.Cells.Item(10,1).Value = "Claudiene"
$ws.Cells.Item(10,3).Value = "FreeSpeech multimedia"
$ws.Cells.Item(10,4).Value = "leitte_anonima@lettersweb.com"

# Row 11 - Zywywz Correia de Campos (no email)
$ws.Cells.Item(11,1).Value = "Zywywz"
$ws.Cells.Item(11,2).Value = "Correia de Campos"
$ws.Cells.Item(11,3).Value = "Z y Z Contabilidade S/A"

# Row 12 - Maria Rosângela Américo (no company)
$ws.Cells.Item(12,1).Value = "Maria Rosângela"
$ws.Cells.Item(12,2).Value = "Américo"
$ws.Cells.Item(12,4).Value = "roseamerico@jinnmail.com"

# Row 13 - (no name) Pereira
$ws.Cells.Item(13,2).Value = "Pereira"
$ws.Cells.Item(13,3).Value = "FunkAgenda Produções"

# Row 14 - Lucas (no surname, no company)
$ws.Cells.Item(14,1).Value = "Lucas"
$ws.Cells.Item(14,4).Value = "lucas_motta_21@ibapo.ru"

# Row 15 - only company
$ws.Cells.Item(15,3).Value = "Restaurante Al-Bait Maghribiy LTDA"

# ---- Hyperlinks on EMAIL column (mailto:) ----
$ws.Hyperlinks.Add($ws.Cells.Item(2,4),  "mailto:tonin@coldmail.froz.en")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4),  "mailto:leitte_anonima@lettersweb.com")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4),  "mailto:zivio.do.agape@jinnmail.comma")
$ws.Hyperlinks.Add($ws.Cells.Item(5,4),  "mailto:roseamerico@jinnmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(6,4),  "mailto:wbeatzz@coldmail.froz.en")
$ws.Hyperlinks.Add($ws.Cells.Item(7,4),  "mailto:lucas_motta_21@ibapo.ru")
$ws.Hyperlinks.Add($ws.Cells.Item(8,4),  "mailto:krmghrib_morroco@jinnmail.comma")
$ws.Hyperlinks.Add($ws.Cells.Item(9,4),  "mailto:bruko_ultra@coldmail.froz.en")
$ws.Hyperlinks.Add($ws.Cells.Item(10,4), "mailto:leitte_anonima@lettersweb.com")
$ws.Hyperlinks.Add($ws.Cells.Item(12,4), "mailto:roseamerico@jinnmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(14,4), "mailto:lucas_motta_21@ibapo.ru")

# ---- Fills ----
$ws.Range("A2:D8").Interior.ThemeColor = 10
$ws.Range("A10:D15").Interior.ThemeColor = 6
$ws.Range("A9:D9").Interior.ThemeColor = 9

# ---- Header bold ----
$ws.Range("A1:D1").Font.Bold = $true

# ---- Column widths (engine quantises to 1/6 character units, so these are
#      the closest achievable approximations of the recorded 17.85546875 /
#      32.42578125 bestFit widths) ----
$ws.Columns.Item(2).ColumnWidth = 17.0
$ws.Columns.Item(3).ColumnWidth = 31.666666666666668

# ---- Selection ----
$ws.Range("D9").Select()
